$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.803267002105713
$ws.Range("B1").Value = 4.641120433807373
$ws.Range("C1").Value = 2.802369356155396
$ws.Range("D1").Value = 1.433313369750977
$ws.Range("E1").Value = 1.056817889213562
